# "south, yuktahar mistakes rectified"
# Correct two menu item names on the AllPages sheet and move the
# active cell selection, matching the authoring app's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllPages")

# Wednesday Breakfast ITEM 1: "Veg. Upma" -> "Veg. Poha"
$ws.Range("D2").Value = "Veg. Poha"

# Friday Breakfast ITEM 2: "Ragi malt" -> "Ragi malt & Daliya"
$ws.Range("G3").Value = "Ragi malt & Daliya"

# Update the last active cell selection saved with the workbook.
$ws.Range("C16").Select()
